$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s301 = @'
Armazenamento de imagens usando Amazon S3
'@

$s302 = @'
Criando um bucket no S3
'@

$s303 = @'
0:31
7. Armazenamento de imagens usando Amazon S3
82. Criando um bucket no S3
criação do bucket no s3 (AWS - Amazon) - pesquisar por all services (todos os serviços) e progurar a opção Armazenamento>S3
'@

$s304 = @'
Setup do IAM - Identity Access Management
'@

$s305 = @'
1:02
7. Armazenamento de imagens usando Amazon S3
83. Setup do IAM - Identity Access Management
configuração do IAM - Identity Access Management do AWS
'@

$s306 = @'

1:22
7. Armazenamento de imagens usando Amazon S3
83. Setup do IAM - Identity Access Management
configurar MFA (Multi Factory Autenticator) - multiplica o grau de segurança do acesso ao AWS ... faz com que o acesso ao AWS nao seja apenas pela senha, mas também por um codigo gerado por um aplicativo no smartphone
'@

$s307 = @'

1:57
7. Armazenamento de imagens usando Amazon S3
83. Setup do IAM - Identity Access Management
necessário um aplicativo de authenticação .. o professor recomenda o uso do google authenticator - app que gera codigos de acesso para inserir na conta AWS ... como se fosse uma authenticação de varios fatores
'@

$s308 = @'
3:03
7. Armazenamento de imagens usando Amazon S3
83. Setup do IAM - Identity Access Management
criação de grupo de usuarios - no caso foi criado um grupo chamado "developers" ... simulando que vamos dar permissao para os desenvolvedores do projeto
'@

$s309 = @'
3:16
7. Armazenamento de imagens usando Amazon S3
83. Setup do IAM - Identity Access Management
adiciona a permissão AmazonS3FullAccess ao grupo developers
'@

$s310 = @'
3:58
7. Armazenamento de imagens usando Amazon S3
83. Setup do IAM - Identity Access Management
criação de usuário - onde damos o nome de "curso-spring-ionic-user" ... marcar a caixa "Acesso Programático", ou seja, alegando que este usuário tera acesso a algum programa que estivermos desenvolvendo
'@

$s311 = @'
6:01
7. Armazenamento de imagens usando Amazon S3
83. Setup do IAM - Identity Access Management
IMPORTANTISSIMO: após criação do usuário, o AWS mostra os dados do novo usuario, a ID da chave de acesso e a chave de acesso secreta ... é importantissimo salvar estes dados em um local seguro pois ele demonstra esses dados somente na criação, ou seja, somente uma vez.
'@

$s312 = @'
6:15
7. Armazenamento de imagens usando Amazon S3
83. Setup do IAM - Identity Access Management
o ultimo item "Aplicar uma politica de senhas do IAM" é opcional ... ele tem a função de adicionar politicas de senha aos usuários, tais como quantidade de caracteres obrigatórios na senha, letra maiúscula ou minuscula, números obrigatórios, etc...
'@

$s313 = @'
2:24
7. Armazenamento de imagens usando Amazon S3
84. Salvando primeiro arquivo no S3
alteração do arquivo application.properties - atenção ao id e access key do AWS ... inserir via variavel de ambiente por segurança... e para encontrar a região que foi inserida no S3 basta acessar o link: http://docs.aws.amazon.com/AWSEC2/latest/UserGuide/using-regions-availability-zones.html (que no caso foi utilizada a região South America (São Paulo) codigo: "sa-east-1"
'@

$s314 = @'
3:56
7. Armazenamento de imagens usando Amazon S3
84. Salvando primeiro arquivo no S3
criação da classe S3Config
'@

$s315 = @'

7:22
7. Armazenamento de imagens usando Amazon S3
84. Salvando primeiro arquivo no S3
criação da classe S3Service
'@

$s316 = @'
9:31
7. Armazenamento de imagens usando Amazon S3
84. Salvando primeiro arquivo no S3
ao implementar o metodo uploadFile, utilizar o file do java.io e no s3client.putObject utilizar o que recebe como parametro um putObject .. pois possui 4 metodos putObject
'@

$s317 = @'
12:37
7. Armazenamento de imagens usando Amazon S3
84. Salvando primeiro arquivo no S3
injeta um teste basico direto na classe principal de projeto CursomcApplication .. apenas para testar os arquivos da aula e verificar se esta fazendo um upload de imagem
'@

$s318 = @'
Salvando primeiro arquivo no S3
'@

$ws.Range("B203:G203").Copy() | Out-Null
$ws.Range("B204:G217").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B204").Value = 7
$ws.Range("C204").Value = $s301
$ws.Range("D204").Value = 82
$ws.Range("E204").Value = $s302
$ws.Range("F204").Value = $s303
$ws.Rows.Item(204).RowHeight = 75

$ws.Range("B205").Value = 7
$ws.Range("C205").Value = $s301
$ws.Range("D205").Value = 83
$ws.Range("E205").Value = $s304
$ws.Range("F205").Value = $s305
$ws.Rows.Item(205).RowHeight = 60

$ws.Range("B206").Value = 7
$ws.Range("C206").Value = $s301
$ws.Range("D206").Value = 83
$ws.Range("E206").Value = $s304
$ws.Range("F206").Value = $s306
$ws.Rows.Item(206).RowHeight = 120

$ws.Range("B207").Value = 7
$ws.Range("C207").Value = $s301
$ws.Range("D207").Value = 83
$ws.Range("E207").Value = $s304
$ws.Range("F207").Value = $s307
$ws.Rows.Item(207).RowHeight = 120

$ws.Range("B208").Value = 7
$ws.Range("C208").Value = $s301
$ws.Range("D208").Value = 83
$ws.Range("E208").Value = $s304
$ws.Range("F208").Value = $s308
$ws.Rows.Item(208).RowHeight = 90

$ws.Range("B209").Value = 7
$ws.Range("C209").Value = $s301
$ws.Range("D209").Value = 83
$ws.Range("E209").Value = $s304
$ws.Range("F209").Value = $s309
$ws.Rows.Item(209).RowHeight = 60

$ws.Range("B210").Value = 7
$ws.Range("C210").Value = $s301
$ws.Range("D210").Value = 83
$ws.Range("E210").Value = $s304
$ws.Range("F210").Value = $s310
$ws.Rows.Item(210).RowHeight = 105

$ws.Range("B211").Value = 7
$ws.Range("C211").Value = $s301
$ws.Range("D211").Value = 83
$ws.Range("E211").Value = $s304
$ws.Range("F211").Value = $s311
$ws.Rows.Item(211).RowHeight = 120

$ws.Range("B212").Value = 7
$ws.Range("C212").Value = $s301
$ws.Range("D212").Value = 83
$ws.Range("E212").Value = $s304
$ws.Range("F212").Value = $s312
$ws.Rows.Item(212).RowHeight = 105

$ws.Range("B213").Value = 7
$ws.Range("C213").Value = $s301
$ws.Range("D213").Value = 84
$ws.Range("F213").Value = $s313
$ws.Rows.Item(213).RowHeight = 150

$ws.Range("B214").Value = 7
$ws.Range("C214").Value = $s301
$ws.Range("D214").Value = 84
$ws.Range("F214").Value = $s314
$ws.Rows.Item(214).RowHeight = 60

$ws.Range("B215").Value = 7
$ws.Range("C215").Value = $s301
$ws.Range("D215").Value = 84
$ws.Range("F215").Value = $s315
$ws.Rows.Item(215).RowHeight = 75

$ws.Range("B216").Value = 7
$ws.Range("C216").Value = $s301
$ws.Range("D216").Value = 84
$ws.Range("F216").Value = $s316
$ws.Rows.Item(216).RowHeight = 90

$ws.Range("B217").Value = 7
$ws.Range("C217").Value = $s301
$ws.Range("D217").Value = 84
$ws.Range("F217").Value = $s317
$ws.Rows.Item(217).RowHeight = 90

$ws.Range("E213:E217").Value = $s318

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B1:G217"))

$ws.Range("E213").Select()
$excel.ActiveWindow.ScrollRow = 211

Write-Host "Done applying edits"